# Daily attendance processing - 2026-01-28 04:49:31
# Swap the order of the "Recorded By" names in column G from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# wherever that exact value appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$old = "dnasr281@gmail.com, System"
$new = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Value2 -eq $old) {
        $cell.Value = $new
    }
}
